$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.04563514837453
$ws.Range("D2").Value = 1.048989521873806
$ws.Range("E2").Value = 1.043272374191925
$ws.Range("F2").Value = 1.05726675085957
$ws.Range("I2").Value = 1.038945101223042
$ws.Range("J2").Value = 1.05069383003779
$ws.Range("K2").Value = 1.051747829631309
$ws.Range("L2").Value = 1.046046707304126
$ws.Range("M2").Value = 1.060002189996293
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.047219679848313
$ws.Range("D3").Value = 1.050208786366259
$ws.Range("E3").Value = 1.044642879433515
$ws.Range("F3").Value = 1.058716748892199
$ws.Range("I3").Value = 1.039344011464254
$ws.Range("J3").Value = 1.051922970927945
$ws.Range("K3").Value = 1.052778198188464
$ws.Range("L3").Value = 1.047226751078704
$ws.Range("M3").Value = 1.061264374463785
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.048243119561685
$ws.Range("D4").Value = 1.050995869573815
$ws.Range("E4").Value = 1.045528195622349
$ws.Range("F4").Value = 1.059653462628366
$ws.Range("I4").Value = 1.039599608590579
$ws.Range("J4").Value = 1.0527161008738
$ws.Range("K4").Value = 1.053442493544437
$ws.Range("L4").Value = 1.04798830091633
$ws.Range("M4").Value = 1.062079045932165
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.048672938301473
$ws.Range("D5").Value = 1.05132631978107
$ws.Range("E5").Value = 1.04590003263182
$ws.Range("F5").Value = 1.060046897896429
$ws.Range("I5").Value = 1.039706459944936
$ws.Range("J5").Value = 1.053049011617624
$ws.Range("K5").Value = 1.053721189366414
$ws.Range("L5").Value = 1.04830798026662
$ws.Range("M5").Value = 1.06242105085341
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.048745081507456
$ws.Range("D6").Value = 1.051381778230288
$ws.Range("E6").Value = 1.045962445395373
$ws.Range("F6").Value = 1.060112936565062
$ws.Range("I6").Value = 1.039724365552408
$ws.Range("J6").Value = 1.053104878445934
$ws.Range("K6").Value = 1.053767950134936
$ws.Range("L6").Value = 1.048361628107236
$ws.Range("M6").Value = 1.062478446863439
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.048248864519481
$ws.Range("D7").Value = 1.051000286785134
$ws.Range("E7").Value = 1.045533165492534
$ws.Range("F7").Value = 1.059658721131713
$ws.Range("I7").Value = 1.039601038704505
$ws.Range("J7").Value = 1.052720551282191
$ws.Range("K7").Value = 1.053446219739639
$ws.Range("L7").Value = 1.047992574350272
$ws.Range("M7").Value = 1.062083617707869
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.046171038684981
$ws.Range("D8").Value = 1.049401967457021
$ws.Range("E8").Value = 1.043735856296065
$ws.Range("F8").Value = 1.057757106253327
$ws.Range("I8").Value = 1.039080439140628
$ws.Range("J8").Value = 1.051109685113581
$ws.Range("K8").Value = 1.052096552556713
$ws.Range("L8").Value = 1.046445930423233
$ws.Range("M8").Value = 1.060429179023087
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.042494985065504
$ws.Range("D9").Value = 1.046570976444088
$ws.Range("E9").Value = 1.040557007290882
$ws.Range("F9").Value = 1.054394144293653
$ws.Range("I9").Value = 1.038143620304707
$ws.Range("J9").Value = 1.048253913603922
$ws.Range("K9").Value = 1.049699451605135
$ws.Range("L9").Value = 1.04370479083482
$ws.Range("M9").Value = 1.057497862877584
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.040033825574569
$ws.Range("D10").Value = 1.044673467546509
$ws.Range("E10").Value = 1.038429402148725
$ws.Range("F10").Value = 1.052143572707678
$ws.Range("I10").Value = 1.037505821568434
$ws.Range("J10").Value = 1.046338036825168
$ws.Range("K10").Value = 1.048088360404429
$ws.Range("L10").Value = 1.041866343152299
$ws.Range("M10").Value = 1.055532466712602
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.038965501575169
$ws.Range("D11").Value = 1.043849320354409
$ws.Range("E11").Value = 1.037506035510072
$ws.Range("F11").Value = 1.051166903744555
$ws.Range("I11").Value = 1.037226465435273
$ws.Range("J11").Value = 1.045505486591482
$ws.Range("K11").Value = 1.047387567017713
$ws.Range("L11").Value = 1.0410675649056
$ws.Range("M11").Value = 1.054678676948245
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.038568272322506
$ws.Range("D12").Value = 1.043542809940865
$ws.Range("E12").Value = 1.037162731520866
$ws.Range("F12").Value = 1.050803792036693
$ws.Range("I12").Value = 1.037122218179207
$ws.Range("J12").Value = 1.045195786394257
$ws.Range("K12").Value = 1.047126776011954
$ws.Range("L12").Value = 1.040770446352184
$ws.Range("M12").Value = 1.054361118131176
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.038653497913347
$ws.Range("D13").Value = 1.043608575034437
$ws.Range("E13").Value = 1.037236386223023
$ws.Range("F13").Value = 1.05088169596752
$ws.Range("I13").Value = 1.037144601434201
$ws.Range("J13").Value = 1.045262238823523
$ws.Range("K13").Value = 1.04718273864014
$ws.Range("L13").Value = 1.040834198284154
$ws.Range("M13").Value = 1.054429254868613
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.038932674834493
$ws.Range("D14").Value = 1.043823992016807
$ws.Range("E14").Value = 1.037477664567828
$ws.Range("F14").Value = 1.051136895672603
$ws.Range("I14").Value = 1.037217858176168
$ws.Range("J14").Value = 1.045479895988771
$ws.Range("K14").Value = 1.047366019891234
$ws.Range("L14").Value = 1.041043013543421
$ws.Range("M14").Value = 1.054652436122184
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.039104630779596
$ws.Range("D15").Value = 1.043956666216998
$ws.Range("E15").Value = 1.037626280846235
$ws.Range("F15").Value = 1.05129408817739
$ws.Range("I15").Value = 1.037262930106415
$ws.Range("J15").Value = 1.045613941330965
$ws.Range("K15").Value = 1.047478881012615
$ws.Range("L15").Value = 1.041171616042078
$ws.Range("M15").Value = 1.054789889110964
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.040104670405722
$ws.Range("D16").Value = 1.044728109860073
$ws.Range("E16").Value = 1.038490637943634
$ws.Range("F16").Value = 1.052208344775701
$ws.Range("I16").Value = 1.037524294117334
$ws.Range("J16").Value = 1.046393227265234
$ws.Range("K16").Value = 1.048134802124072
$ws.Range("L16").Value = 1.041919297451355
$ws.Range("M16").Value = 1.055589071040923
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.040731257661335
$ws.Range("D17").Value = 1.045211337668915
$ws.Range("E17").Value = 1.039032257821976
$ws.Range("F17").Value = 1.052781249816369
$ws.Range("I17").Value = 1.037687385712982
$ws.Range("J17").Value = 1.046881253275848
$ws.Range("K17").Value = 1.048545387500966
$ws.Range("L17").Value = 1.042387564611167
$ws.Range("M17").Value = 1.056089631716147
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.04109648285863
$ws.Range("D18").Value = 1.045492954310586
$ws.Range("E18").Value = 1.039347973572571
$ws.Range("F18").Value = 1.053115208189174
$ws.Range("I18").Value = 1.037782207202127
$ws.Range("J18").Value = 1.047165625310369
$ws.Range("K18").Value = 1.04878456865529
$ws.Range("L18").Value = 1.042660435218413
$ws.Range("M18").Value = 1.056381334549296
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.041220972727472
$ws.Range("D19").Value = 1.045588937502713
$ws.Range("E19").Value = 1.039455590501554
$ws.Range("F19").Value = 1.053229044459483
$ws.Range("I19").Value = 1.037814486896561
$ws.Range("J19").Value = 1.047262540702198
$ws.Range("K19").Value = 1.048866071506495
$ws.Range("L19").Value = 1.042753432905748
$ws.Range("M19").Value = 1.056480752861326
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.040664056985427
$ws.Range("D20").Value = 1.045159516982542
$ws.Range("E20").Value = 1.038974168090296
$ws.Range("F20").Value = 1.052719804052167
$ws.Range("I20").Value = 1.037669919315727
$ws.Range("J20").Value = 1.046828922275315
$ws.Range("K20").Value = 1.048501367316059
$ws.Range("L20").Value = 1.042337351106391
$ws.Range("M20").Value = 1.056035953841832
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.038850475495543
$ws.Range("D21").Value = 1.043760567788421
$ws.Range("E21").Value = 1.037406623140203
$ws.Range("F21").Value = 1.051061754998369
$ws.Range("I21").Value = 1.037196299223431
$ws.Range("J21").Value = 1.045415813981646
$ws.Range("K21").Value = 1.047312061582468
$ws.Range("L21").Value = 1.040981534225902
$ws.Range("M21").Value = 1.054586726559712
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.037707849105031
$ws.Range("D22").Value = 1.042878758857782
$ws.Range("E22").Value = 1.036419163366944
$ws.Range("F22").Value = 1.050017340280795
$ws.Range("I22").Value = 1.036895725334611
$ws.Range("J22").Value = 1.044524705745158
$ws.Range("K22").Value = 1.046561487311232
$ws.Range("L22").Value = 1.040126662976729
$ws.Range("M22").Value = 1.053673086547199
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.03831380419213
$ws.Range("D23").Value = 1.043346436859717
$ws.Range("E23").Value = 1.036942815921071
$ws.Range("F23").Value = 1.050571190701788
$ws.Range("I23").Value = 1.037055330867217
$ws.Range("J23").Value = 1.044997351432845
$ws.Range("K23").Value = 1.046959649810418
$ws.Range("L23").Value = 1.040580078158641
$ws.Range("M23").Value = 1.054157659836257
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.040694422858712
$ws.Range("D24").Value = 1.045182933267957
$ws.Range("E24").Value = 1.03900041696381
$ws.Range("F24").Value = 1.052747569389288
$ws.Range("I24").Value = 1.03767781257775
$ws.Range("J24").Value = 1.046852569283527
$ws.Range("K24").Value = 1.048521259089533
$ws.Range("L24").Value = 1.042360041236567
$ws.Range("M24").Value = 1.056060209383981
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.043447131796739
$ws.Range("D25").Value = 1.047304620475771
$ws.Range("E25").Value = 1.041380257306682
$ws.Range("F25").Value = 1.05526503043806
$ws.Range("I25").Value = 1.038388133621230
$ws.Range("J25").Value = 1.048994285846524
$ws.Range("K25").Value = 1.050321427589336
$ws.Range("L25").Value = 1.044415352081281
$ws.Range("M25").Value = 1.058257617455097
